# "show update table by switch column"
#
# A new product row ("Fahrrad" / qty=1 / price=20.0 / category="Keller")
# is inserted as the new row 4 of the product table, pushing the
# existing rows 4-8 (Fressen, Halsband, Katzenspielzeug, Maus, Toilette)
# down to rows 5-9, and the two summary rows down from 9-10 to 10-11.
# The summary counters are recalculated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4, shifting rows 4.. down by one.
$ws.Rows("4:4").Insert()

# The table alternates row styling every other row (zebra stripes).
# Because every data row below the insertion point moved down by one,
# the stripe pattern needs to be re-applied from row 4 through row 9
# using the two already-present formatting templates (row 2 / row 3).
$ws.Range("A2:F2").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)
$ws.Range("A6:F6").PasteSpecial(-4122)
$ws.Range("A8:F8").PasteSpecial(-4122)

$ws.Range("A3:F3").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)
$ws.Range("A7:F7").PasteSpecial(-4122)
$ws.Range("A9:F9").PasteSpecial(-4122)

# Fill in the new product row.
$ws.Range("A4").Value = "Fahrrad"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 20.0
$ws.Range("E4").Value = "Keller"

# Restore the row heights used throughout the table.
$ws.Rows("4:9").RowHeight = 20
$ws.Rows("10:11").RowHeight = 25

# Update the "Anzahl der Produkte:" / "Gesamtwert:" summary values to
# reflect the newly added row.
$ws.Range("F10").Value = 42
$ws.Range("F11").Value = 359
